# Apply updated crypto price/volume figures per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.136.94"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "2.223.23"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "293.95"
$ws.Range("E5").Value = "  +1.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.86"
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.513"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.67"
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "50.86"
$ws.Range("E11").Value = "  +6.92%  "
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("E13").Value = "  +3.97%  "
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "2.565.63"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.83"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "2.211.18"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").Value = "40.067.08"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.27"
$ws.Range("E21").Value = "  -2.27%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.68"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.92"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.48"
$ws.Range("E26").Value = "  +1.91%  "
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.18"
$ws.Range("E28").Value = "  +3.49%  "
$ws.Range("E29").Value = "  +1.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.06"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.27"
$ws.Range("E31").Value = "  +4.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.85"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.96"
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.05"
$ws.Range("E35").Value = "  +8.61%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.77"
$ws.Range("E39").Value = "  +4.75%  "
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "2.076.11"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.73"
$ws.Range("E44").Value = "  +13.56%  "
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.99"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.75"
$ws.Range("E47").Value = "  +4.17%  "
$ws.Range("E48").Value = "  -11.09%  "
$ws.Range("D49").Value = "2.437.28"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("E50").Value = "  +5.30%  "
$ws.Range("E51").Value = "  +2.10%  "
